$d = $word.ActiveDocument

# 1) Append a new empty paragraph right after the "Dia 17/09" line.
$d.Content.Find.Execute("Dia 17/09: 1hr (1 dia)", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Dia 17/09: 1hr (1 dia)^p", 2)

# 2) Bump the hour count on that line from 1hr to 2hr, splitting the run the
#    same way Word does when a user selects just the digit and retypes it
#    (toggling a character property forces a fresh run instead of merging
#    back into its neighbours).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Dia 17/09: 1hr (1 dia)*") {
        $pStart = $p.Range.Start
        $digit = $d.Range($pStart + 11, $pStart + 12)
        $digit.Bold = 1
        $digit.Text = "2"
        $digit.Bold = 0
    }
}
